# Add "DateDébut" / "DateFin" columns (M, N) and three new student rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a plain text value into a cell without Excel re-interpreting
# it as a date/number and without leaving a non-default cell style behind.
function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# --- New columns: headers ---------------------------------------------------
$ws.Range("M1").Value = "DateDébut"
$ws.Range("N1").Value = "DateFin"

# --- Insert a new row for the new first student (ABBE Tristan) -------------
# This shifts the existing "MONTBULEAU--GENTELET" row from row 2 to row 3,
# preserving all of its existing cell values.
$ws.Rows(2).Insert()

# Row 2 : ABBE Tristan
$ws.Range("A2").Value = "ABBE"
$ws.Range("B2").Value = "TRISTAN"
$ws.Range("C2").Value = 2025
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = "GPhy"
$ws.Range("G2").Value = "stage"
$ws.Range("H2").Value = "SANOFI"
$ws.Range("I2").Value = "Gentilly (94)"
$ws.Range("J2").Value = "BELLOCQ"
$ws.Range("L2").Value = "GENIET"
Set-TextValue "M2" "2024-05-21"
Set-TextValue "N2" "2024-08-23"

# Row 3 : MONTBULEAU--GENTELET Titouan (shifted down from row 2).
# Nature changes from "stage" to "apprentissage"; add the new date columns.
$ws.Range("G3").Value = "apprentissage"
Set-TextValue "M3" "2024-05-21"
Set-TextValue "N3" "2024-08-31"

# Row 4 : SEVILLA Mathieu
$ws.Range("A4").Value = "SEVILLA"
$ws.Range("B4").Value = "MATHIEU"
$ws.Range("C4").Value = 2025
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = "GPhy"
$ws.Range("G4").Value = "stage"
$ws.Range("H4").Value = "LABVANTAGE"
$ws.Range("I4").Value = "Paris"
$ws.Range("J4").Value = "MAIRE"
Set-TextValue "K4" "+(33)6 74.01.53.57"
$ws.Range("L4").Value = "URRUTY"
Set-TextValue "M4" "2024-05-21"
Set-TextValue "N4" "2024-08-23"

# Row 5 : SOCHARD Ophelie
$ws.Range("A5").Value = "SOCHARD"
$ws.Range("B5").Value = "OPHELIE"
$ws.Range("C5").Value = 2025
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = "GPhy"
$ws.Range("G5").Value = "stage"
$ws.Range("H5").Value = "AMEXIO"
$ws.Range("I5").Value = "Paris, Montpellier, Lille et Nantes"
$ws.Range("J5").Value = "SOCHARD"
Set-TextValue "K5" "01 81 69 86 00"
$ws.Range("L5").Value = "GENIET"
Set-TextValue "M5" "2024-05-27"
Set-TextValue "N5" "2024-08-16"
